# Add extent reporting logic along with few more cases
# - Adds two new worksheets: BOOKFLIGHT and TESTRUNNER (with their data)
# - Updates the active sheet / selections to reflect the new TESTRUNNER tab
# - Sets TESTRUNNER's column widths

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- BOOKFLIGHT sheet -------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "BOOKFLIGHT"

$ws2.Range("A1").Value = "departcity"
$ws2.Range("B1").Value = "arrivalcity"
$ws2.Range("A2").Value = "Patna"
$ws2.Range("B2").Value = "Delhi"
$ws2.Range("A3").Value = "Mumbai"
$ws2.Range("B3").Value = "Pune"
$ws2.Range("A4").Value = "Kolkata"
$ws2.Range("B4").Value = "Lucknow"

[void]$ws2.Range("G11").Select()

# --- TESTRUNNER sheet --------------------------------------------------
$ws3 = $wb.Worksheets.Add([Type]::Missing, $ws2)
$ws3.Name = "TESTRUNNER"

$ws3.Range("A1").Value = "testName"
$ws3.Range("B1").Value = "browser"
$ws3.Range("A2").Value = "bookFlight"
$ws3.Range("B2").Value = "chrome"
$ws3.Range("A3").Value = "addToCart"
$ws3.Range("B3").Value = "edge"

$ws3.Columns.Item(1).ColumnWidth = 41.833333333333336
$ws3.Columns.Item(2).ColumnWidth = 21.5

# --- Selections / active sheet -----------------------------------------
# Sheet1 keeps a new selection but is no longer the active tab.
[void]$ws1.Range("A6:B8").Select()

# TESTRUNNER becomes the active (selected) tab, with its own selection.
[void]$ws3.Activate()
[void]$ws3.Range("A9").Select()
